$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit cyclically rotates the data that lives in rows 3, 4, 6 and 7
# (row 5 is untouched):
#   new row3 <- old row7
#   new row4 <- old row6
#   new row6 <- old row3
#   new row7 <- old row4
#
# Columns that actually carry data for these rows: A, B, D, E, F, G, H, M, Q, R
# (M only has a value in the original row 6 - "färska gnagspår" - and ends up
# in row 4 after the rotation).

$cols = @("A", "B", "D", "E", "F", "G", "H", "M", "Q", "R")

# Snapshot every source value BEFORE any writes happen, so the rotation
# doesn't clobber a value before it has been read.
$old3 = @{}
$old4 = @{}
$old6 = @{}
$old7 = @{}
foreach ($col in $cols) {
    $old3[$col] = $ws.Range($col + "3").Value2
    $old4[$col] = $ws.Range($col + "4").Value2
    $old6[$col] = $ws.Range($col + "6").Value2
    $old7[$col] = $ws.Range($col + "7").Value2
}

foreach ($col in $cols) {
    $addr3 = $col + "3"
    $addr4 = $col + "4"
    $addr6 = $col + "6"
    $addr7 = $col + "7"

    if ($col -eq "M") {
        # Column M only ever has a value for one row at a time - use
        # ClearContents for the destinations that should end up empty.
        if ($null -eq $old7[$col]) { $ws.Range($addr3).ClearContents() } else { $ws.Range($addr3).Value2 = $old7[$col] }
        if ($null -eq $old6[$col]) { $ws.Range($addr4).ClearContents() } else { $ws.Range($addr4).Value2 = $old6[$col] }
        if ($null -eq $old3[$col]) { $ws.Range($addr6).ClearContents() } else { $ws.Range($addr6).Value2 = $old3[$col] }
        if ($null -eq $old4[$col]) { $ws.Range($addr7).ClearContents() } else { $ws.Range($addr7).Value2 = $old4[$col] }
    } else {
        $ws.Range($addr3).Value2 = $old7[$col]
        $ws.Range($addr4).Value2 = $old6[$col]
        $ws.Range($addr6).Value2 = $old3[$col]
        $ws.Range($addr7).Value2 = $old4[$col]
    }
}
